# Apply "first values to Singapore archetypes database" edit
$wb = $excel.ActiveWorkbook

# --- ARCHITECTURE sheet ---
$ws1 = $wb.Worksheets.Item("ARCHITECTURE")

# Row -> building_use name (for reference / readability only)
# 2 MULTI_RES, 3 SINGLE_RES, 4 HOTEL, 5 OFFICE, 6 RETAIL, 7 FOODSTORE,
# 8 RESTAURANT, 9 INDUSTRIAL, 10 SCHOOL, 11 HOSPITAL, 12 GYM, 13 SWIMMING,
# 14 SERVERROOM, 15 PARKING, 16 COOLROOM, 17 LAB, 18 MUSEUM, 19 LIBRARY

# Column G = type_cons -> "T2" for every data row
for ($r = 2; $r -le 19; $r++) {
    $ws1.Range("G$r").Value = "T2"
}

# Column J = type_roof -> "T7" for every data row
for ($r = 2; $r -le 19; $r++) {
    $ws1.Range("J$r").Value = "T7"
}

# Column K = type_wall -> "T2" for rows 2-3, "T6" for several other rows
$ws1.Range("K2").Value = "T2"
$ws1.Range("K3").Value = "T2"
$ws1.Range("K5").Value = "T6"
$ws1.Range("K6").Value = "T6"
$ws1.Range("K8").Value = "T6"
$ws1.Range("K10").Value = "T6"
$ws1.Range("K12").Value = "T6"
$ws1.Range("K17").Value = "T6"
$ws1.Range("K18").Value = "T6"
$ws1.Range("K19").Value = "T6"

# Row 5 (OFFICE) extra values
$ws1.Range("F5").Value = 0.35
$ws1.Range("H5").Value = "T1"
$ws1.Range("I5").Value = "T2"

# Row 17 (LAB) extra value
$ws1.Range("F17").Value = 0.11

# --- HVAC sheet ---
$ws2 = $wb.Worksheets.Item("HVAC")

# Column E = type_hs -> "T0" for every data row
for ($r = 2; $r -le 19; $r++) {
    $ws2.Range("E$r").Value = "T0"
}

# Update the active selection on this sheet (was A2:A19, now E2:E19)
$ws2.Activate()
$ws2.Range("E2:E19").Select()

# ARCHITECTURE stays the selected/visible tab; update its selection last
# (was G23, now G19)
$ws1.Activate()
$ws1.Range("G19").Select()

$wb.Save()
